# maglev_results.xlsx -- complete dynamics linearization
# Extends the y10/y20 operating-point sweep from +/-1.4 out to +/-2.2 and
# refreshes every derived column (cond(o), cond(c), TF12_num, TF21_num,
# tau_dom, dom eig) for the full table, rows 2-23, columns A:R.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Column width adjustments (A widened; I:J and N:O spacer columns widened; R now matches Q) ---
$ws.Columns.Item(1).ColumnWidth  = 9.666666666666666
$ws.Columns.Item(9).ColumnWidth  = 14.833333333333334
$ws.Columns.Item(10).ColumnWidth = 14.833333333333334
$ws.Columns.Item(14).ColumnWidth = 14.833333333333334
$ws.Columns.Item(15).ColumnWidth = 14.833333333333334
$ws.Columns.Item(18).ColumnWidth = 11.666666666666666

# --- Results table: rows 2-23, columns A:R ---
# A=y10  B=y20  C=u10  D=u20  E=cond(o)  F=cond(c)
# G:K=TF12_num (G,H spacer; I,J,K coefficients)
# L:P=TF21_num (L,M spacer; N,O,P coefficients)
# Q=tau_dom  R=dom eig
$results = New-Object 'object[,]' 22,18

# Row 2  y10=0.1  y20=-0.1
$results[0,0] = 0.10000000000000001
$results[0,1] = -0.10000000000000001
$results[0,2] = 4005.9683493065641
$results[0,3] = 2204.8371461562447
$results[0,4] = 3.8041409678215623
$results[0,5] = 2.9861922376603327
$results[0,6] = 0
$results[0,7] = 0
$results[0,8] = -0.0000470565962669699399843
$results[0,9] = -0.0000117641490667424426444
$results[0,10] = 0.0025186071013528696
$results[0,11] = 0
$results[0,12] = 0
$results[0,13] = -0.0000470565962669699332081
$results[0,14] = -0.0000117641490667423155895
$results[0,15] = 0.0019267237441451676
$results[0,16] = 0.55019652757935256
$results[0,17] = 1.8175323722954142

# Row 3  y10=0.2  y20=-0.2
$results[1,0] = 0.20000000000000001
$results[1,1] = -0.20000000000000001
$results[1,2] = 4318.9229147504038
$results[1,3] = 2304.4777323350422
$results[1,4] = 4.0128383554782312
$results[1,5] = 3.1519516220210639
$results[1,6] = 0
$results[1,7] = 0
$results[1,8] = -0.0000481110449103133296523
$results[1,9] = -0.0000120277612275782477098
$results[1,10] = 0.002512186949154884
$results[1,11] = 0
$results[1,12] = 0
$results[1,13] = -0.0000481110449103133296523
$results[1,14] = -0.0000120277612275784611621
$results[1,15] = 0.0019160502731082887
$results[1,16] = 0.56628067135680671
$results[1,17] = 1.7659087632357346

# Row 4  y10=0.3  y20=-0.3
$results[2,0] = 0.30000000000000004
$results[2,1] = -0.30000000000000004
$results[2,2] = 4654.6446479958058
$results[2,3] = 2402.6038327594197
$results[2,4] = 4.2526809801803411
$results[2,5] = 3.3427340276185311
$results[2,6] = 0
$results[2,7] = 0
$results[2,8] = -0.0000491951949379011068005
$results[2,9] = -0.0000122987987344751886087
$results[2,10] = 0.0025093798595807594
$results[2,11] = 0
$results[2,12] = 0
$results[2,13] = -0.0000491951949379011000242
$results[2,14] = -0.000012298798734475144563
$results[2,15] = 0.0019087107441542424
$results[2,16] = 0.58414013555999023
$results[2,17] = 1.711917978451083

# Row 5  y10=0.4  y20=-0.4
$results[3,0] = 0.40000000000000002
$results[3,1] = -0.40000000000000002
$results[3,2] = 5014.9519254936731
$results[3,3] = 2498.2889636882178
$results[3,4] = 4.5303576590027888
$results[3,5] = 3.5639702362468553
$results[3,6] = 0
$results[3,7] = 0
$results[3,8] = -0.0000503100559727481558519
$results[3,9] = -0.0000125775139931872625797
$results[3,10] = 0.0025101163174238278
$results[3,11] = 0
$results[3,12] = 0
$results[3,13] = -0.0000503100559727481897333
$results[3,14] = -0.000012577513993187091479
$results[3,15] = 0.0019046321291226521
$results[3,16] = 0.60416191509731698
$results[3,17] = 1.6551854312745322

# Row 6  y10=0.5  y20=-0.5
$results[4,0] = 0.5
$results[4,1] = -0.5
$results[4,2] = 5401.8450121621154
$results[4,3] = 2590.4706060773165
$results[4,4] = 4.854455453582343
$results[4,5] = 3.8226455453650194
$results[4,6] = 0
$results[4,7] = 0
$results[4,8] = -0.0000514566779051693401264
$results[4,9] = -0.0000128641694762922892918
$results[4,10] = 0.0025143461243464942
$results[4,11] = 0
$results[4,12] = 0
$results[4,13] = -0.0000514566779051693401264
$results[4,14] = -0.0000128641694762925179907
$results[4,15] = 0.0019037601375470937
$results[4,16] = 0.62685350612646906
$results[4,17] = 1.5952690544547226

# Row 7  y10=0.6  y20=-0.6
$results[5,0] = 0.59999999999999998
$results[5,1] = -0.59999999999999998
$results[5,2] = 5817.5275845997367
$results[5,3] = 2677.9331233542935
$results[5,4] = 5.2361097865164723
$results[5,5] = 4.1278366215698696
$results[5,6] = 0
$results[5,7] = 0
$results[5,8] = -0.0000526361527387114310167
$results[5,9] = -0.0000131590381846779509278
$results[5,10] = 0.002522037418193737
$results[5,11] = 0
$results[5,12] = 0
$results[5,13] = -0.000052636152738711464898
$results[5,14] = -0.0000131590381846779136583
$results[5,15] = 0.0019060582344581625
$results[5,16] = 0.65289523691174911
$results[5,17] = 1.5316392944296644

# Row 8  y10=0.7  y20=-0.7
$results[6,0] = 0.70000000000000007
$results[6,1] = -0.70000000000000007
$results[6,2] = 6264.4310235282646
$results[6,3] = 2759.2884620353675
$results[6,4] = 5.6898873681970654
$results[6,5] = 4.4914433120504302
$results[6,6] = 0
$results[6,7] = 0
$results[6,8] = -0.0000538496165318323472557
$results[6,9] = -0.0000134624041329579919462
$results[6,10] = 0.0025331758847272254
$results[6,11] = 0
$results[6,12] = 0
$results[6,13] = -0.0000538496165318322591643
$results[6,14] = -0.0000134624041329581122249
$results[6,15] = 0.0019115068480491776
$results[6,16] = 0.68322317208738581
$results[6,17] = 1.46365059156995

# Row 9  y10=0.8  y20=-0.8
$results[7,0] = 0.80000000000000004
$results[7,1] = -0.80000000000000004
$results[7,2] = 6745.241859153075
$results[7,3] = 2832.9543264556473
$results[7,4] = 6.234948507241719
$results[7,5] = 4.9291576827223498
$results[7,6] = 0
$results[7,7] = 0
$results[7,8] = -0.0000550982514408754041617
$results[7,9] = -0.0000137745628602189492962
$results[7,10] = 0.0025477641409715518
$results[7,11] = 0
$results[7,12] = 0
$results[7,13] = -0.0000550982514408754515955
$results[7,14] = -0.0000137745628602189611547
$results[7,15] = 0.0019201027464762796
$results[7,16] = 0.71916612461524143
$results[7,17] = 1.3904993099264882

# Row 10  y10=0.9  y20=-0.9
$results[8,0] = 0.90000000000000002
$results[8,1] = -0.90000000000000002
$results[8,2] = 7262.9328111705154
$results[8,3] = 2897.1294724868117
$results[8,4] = 6.8964539136457619
$results[8,5] = 5.4616504066088929
$results[8,6] = 0
$results[8,7] = 0
$results[8,8] = -0.0000563832878702484437063
$results[8,9] = -0.0000140958219675622616984
$results[8,10] = 0.0025658212720185893
$results[8,11] = 0
$results[8,12] = 0
$results[8,13] = -0.0000563832878702485114689
$results[8,14] = -0.0000140958219675621278672
$results[8,15] = 0.0019318585656977647
$results[8,16] = 0.76268394774560178
$results[8,17] = 1.3111591019528794

# Row 11  y10=1  y20=-1
$results[9,0] = 1
$results[9,1] = -1
$results[9,2] = 7820.7979328741312
$results[9,3] = 2949.7657111225972
$results[9,4] = 7.7068738462980511
$results[9,5] = 6.1157135705857488
$results[9,6] = 0
$results[9,7] = 0
$results[9,8] = -0.0000577060067360906953516
$results[9,9] = -0.0000144265016840225704999
$results[9,10] = 0.0025873825051450231
$results[9,11] = 0
$results[9,12] = 0
$results[9,13] = -0.0000577060067360907698905
$results[9,14] = -0.0000144265016840226416506
$results[9,15] = 0.0019468024722469435
$results[9,16] = 0.816808901272147
$results[9,17] = 1.2242765700062037

# Row 12  y10=1.1  y20=-1.1
$results[10,0] = 1.1000000000000001
$results[10,1] = -1.1000000000000001
$results[10,2] = 8422.4924480112459
$results[10,3] = 2988.5361498564334
$results[10,4] = 8.7058865668479015
$results[10,5] = 6.9243296968215331
$results[10,6] = 0
$results[10,7] = 0
$results[10,8] = -0.0000590677418501249723668
$results[10,9] = -0.0000147669354625313481238
$results[10,10] = 0.0026124990064920711
$results[10,11] = 0
$results[10,12] = 0
$results[10,13] = -0.0000590677418501249181567
$results[10,14] = -0.0000147669354625313345713
$results[10,15] = 0.0019649779462084498
$results[10,16] = 0.88653102777918047
$results[10,17] = 1.1279921048054764

# Row 13  y10=1.2  y20=-1.2
$results[11,0] = 1.2000000000000002
$results[11,1] = -1.2000000000000002
$results[11,2] = 9072.0779618754823
$results[11,3] = 3010.7991262202277
$results[11,4] = 9.9347804799948989
$results[11,5] = 7.9224323678235313
$results[11,6] = 0
$results[11,7] = 0
$results[11,8] = -0.0000604698824308263943356
$results[11,9] = -0.0000151174706077067595202
$results[11,10] = 0.0026412377863381463
$results[11,11] = 0
$results[11,12] = 0
$results[11,13] = -0.0000604698824308264282169
$results[11,14] = -0.0000151174706077066070542
$results[11,15] = 0.0019864436704326869
$results[11,16] = 0.98077036243025495
$results[11,17] = 1.0196066666636376

# Row 14  y10=1.3  y20=-1.3
$results[12,0] = 1.3000000000000003
$results[12,1] = -1.3000000000000003
$results[12,2] = 9774.0738372021879
$results[12,3] = 3013.5572017344052
$results[12,4] = 11.414391132797533
$results[12,5] = 9.1296590311547909
$results[12,6] = 0
$results[12,7] = 0
$results[12,8] = -0.0000619138757495079201888
$results[12,9] = -0.0000154784689373770342573
$results[12,10] = 0.0026736816991370666
$results[12,11] = 0
$results[12,12] = 0
$results[12,13] = -0.0000619138757495079472938
$results[12,14] = -0.0000154784689373771528419
$results[12,15] = 0.0020112735121486777
$results[12,16] = 1.117454321393577
$results[12,17] = 0.89489116544191294

# Row 15  y10=1.4  y20=-1.4
$results[13,0] = 1.4000000000000001
$results[13,1] = -1.4000000000000001
$results[13,2] = 10533.515653927709
$results[13,3] = 2993.4104834709942
$results[13,4] = 13.083928524060575
$results[13,5] = 10.502255888113682
$results[13,6] = 0
$results[13,7] = 0
$results[13,8] = -0.0000634012299194245212454
$results[13,9] = -0.0000158503074798560727131
$results[13,10] = 0.002709929523936274
$results[13,11] = 0
$results[13,12] = 0
$results[13,13] = -0.0000634012299194246567707
$results[13,14] = -0.0000158503074798561065944
$results[13,15] = 0.0020395565825618812
$results[13,16] = 1.3397443477165478
$results[13,17] = 0.74641106096427579

# Row 16  y10=1.5  y20=-1.5
$results[14,0] = 1.5000000000000002
$results[14,1] = -1.5000000000000002
$results[14,2] = 11356.021823632545
$results[14,3] = 2946.503421598944
$results[14,4] = 14.684966376713492
$results[14,5] = 11.842983718706254
$results[14,6] = 0
$results[14,7] = 0
$results[14,8] = -0.0000649335168365375937959
$results[14,9] = -0.0000162333792091345136455
$results[14,10] = 0.0027500961094330824
$results[14,11] = 0
$results[14,12] = 0
$results[14,13] = -0.0000649335168365375395858
$results[14,14] = -0.0000162333792091341545035
$results[14,15] = 0.0020713973586521579
$results[14,16] = 1.7912550861998784
$results[14,17] = 0.55826777978422126

# Row 17  y10=1.6  y20=-1.6
$results[15,0] = 1.6000000000000001
$results[15,1] = -1.6000000000000001
$results[15,2] = 12247.869609212974
$results[15,3] = 2868.4640912095715
$results[15,4] = 15.69566566961147
$results[15,5] = 12.768417751247265
$results[15,6] = 0
$results[15,7] = 0
$results[15,8] = -0.0000665123752811578911222
$results[15,9] = -0.0000166280938202894151823
$results[15,10] = 0.0027943125656287727
$results[15,11] = 0
$results[15,12] = 0
$results[15,13] = -0.0000665123752811583925657
$results[15,14] = -0.000016628093820289716726
$results[15,15] = 0.00210691584907398
$results[15,16] = 3.606090610469928
$results[15,17] = 0.27730861700940035

# Row 18  y10=1.7  y20=-1.7
$results[16,0] = 1.7000000000000002
$results[16,1] = -1.7000000000000002
$results[16,2] = 13216.08201377174
$results[16,3] = 2754.3348012184429
$results[16,4] = 15.627632693193997
$results[16,5] = 13.209996628822301
$results[16,6] = 0
$results[16,7] = 0
$results[16,8] = -0.0000681395141903051271108
$results[16,9] = -0.0000170348785475762207913
$results[16,10] = 0.0028427264806007752
$results[16,11] = 0
$results[16,12] = 0
$results[16,13] = -0.0000681395141903051948734
$results[16,14] = -0.0000170348785475762987184
$results[16,15] = 0.0021462477826069367
$results[16,16] = 8.0000000000000036
$results[16,17] = -0.12499999999999996

# Row 19  y10=1.8  y20=-1.8
$results[17,0] = 1.8000000000000003
$results[17,1] = -1.8000000000000003
$results[17,2] = 14268.527256958694
$results[17,3] = 2598.4926771113714
$results[17,4] = 14.578801077706167
$results[17,5] = 12.656981686379355
$results[17,6] = 0
$results[17,7] = 0
$results[17,8] = -0.0000698167161112895538711
$results[17,9] = -0.0000174541790278224494542
$results[17,10] = 0.0028955021360626359
$results[17,11] = 0
$results[17,12] = 0
$results[17,13] = -0.0000698167161112895403186
$results[17,14] = -0.0000174541790278225714269
$results[17,15] = 0.0021895447927393406
$results[17,16] = 8.0000000000000107
$results[17,17] = -0.12499999999999983

# Row 20  y10=1.9  y20=-1.9
$results[18,0] = 1.9000000000000001
$results[18,1] = -1.9000000000000001
$results[18,2] = 15414.032860746493
$results[18,3] = 2394.5586315304267
$results[18,4] = 13.106170090370204
$results[18,5] = 11.211026461110826
$results[18,6] = 0
$results[18,7] = 0
$results[18,8] = -0.0000715458408477352236863
$results[18,9] = -0.0000178864602119338059216
$results[18,10] = 0.0029528206887606693
$results[18,11] = 0
$results[18,12] = 0
$results[18,13] = -0.0000715458408477351965812
$results[18,14] = -0.0000178864602119338635198
$results[18,15] = 0.002236974565332151
$results[18,16] = 8
$results[18,17] = -0.125

# Row 21  y10=2  y20=-2
$results[19,0] = 2
$results[19,1] = -2
$results[19,2] = 16662.516730687603
$results[19,3] = 2135.2928595825151
$results[19,4] = 11.63792304657953
$results[19,5] = 9.8584694850896213
$results[19,6] = 0
$results[19,7] = 0
$results[19,8] = -0.0000733288293100340700251
$results[19,9] = -0.0000183322073275084531318
$results[19,10] = 0.0030148802758904283
$results[19,11] = 0
$results[19,12] = 0
$results[19,13] = -0.0000733288293100342597605
$results[19,14] = -0.0000183322073275085649401
$results[19,15] = 0.0022887209074297274
$results[19,16] = 7.9999999999999964
$results[19,17] = -0.12500000000000006

# Row 22  y10=2.1  y20=-2.1
$results[20,0] = 2.1000000000000001
$results[20,1] = -2.1000000000000001
$results[20,2] = 18025.138056472686
$results[20,3] = 1812.4746649874137
$results[20,4] = 10.358828073365133
$results[20,5] = 8.8329017791809488
$results[20,6] = 0
$results[20,7] = 0
$results[20,8] = -0.0000751677075830452099325
$results[20,9] = -0.0000187919268957613024831
$results[20,10] = 0.0030818959909812557
$results[20,11] = 0
$results[20,12] = 0
$results[20,13] = -0.0000751677075830453590103
$results[20,14] = -0.0000187919268957615396524
$results[20,15] = 0.0023449836835334119
$results[20,16] = 8
$results[20,17] = -0.125

# Row 23  y10=2.2  y20=-2.2
$results[21,0] = 2.2000000000000002
$results[21,1] = -2.2000000000000002
$results[21,2] = 19514.47138381589
$results[21,3] = 1416.7640273239806
$results[21,4] = 9.3077870157774782
$results[21,5] = 8.0211324992448176
$results[21,6] = 0
$results[21,7] = 0
$results[21,8] = -0.0000770645912247475373081
$results[21,9] = -0.000019266147806186677651
$results[21,10] = 0.0031540996612727185
$results[21,11] = 0
$results[21,12] = 0
$results[21,13] = -0.0000770645912247474695455
$results[21,14] = -0.000019266147806186935149
$results[21,15] = 0.0024059785502134777
$results[21,16] = 8.0000000000000071
$results[21,17] = -0.12499999999999989

$ws.Range("A2:R23").Value = $results

# --- Active selection moves to E10 ---
$ws.Range("E10").Select() | Out-Null
